$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Municipio zone labels) - extend/update to match new layout
$ws.Range("F2").Value = "CAM"
$ws.Range("K2").Value = "centro"
$ws.Range("L2").Value = "centro"
$ws.Range("M2").Value = "centro"

# Row 3 (Direccion) - update F3 and add new addresses H3:M3
$ws.Range("F3").Value = "calle 4 #12-17"
$ws.Range("H3").Value = "calle 13n #18-34"
$ws.Range("I3").Value = "calle 13n #18-35"
$ws.Range("J3").Value = "calle 13n #18-36"
$ws.Range("K3").Value = "calle 13n #18-37"
$ws.Range("L3").Value = "calle 13n #18-38"
$ws.Range("M3").Value = "calle 13n #18-39"

# Update the active selection to match the author's final cursor position
[void]$ws.Range("L7").Select()
